$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the xpath-style selectors to plain CSS id selectors (page object cleanup)
$ws.Range("C1").Value = "#video-0"
$ws.Range("C2").Value = "#video-1"

# Drop the now-unused password column (D) entirely, shifting remaining
# cells left and shrinking the used range from A1:D5 to A1:C5
$ws.Columns.Item(4).Delete()
